# Auto-generated edit script applying scheduled profit-sheet data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 43978.5
$ws.Range("J110").Value = 43978.5
$ws.Range("L110").Value = 43978.5
$ws.Range("N110").Value = -52158.5
$ws.Range("H141").Value = 5427.727
$ws.Range("I141").Value = 5680
$ws.Range("K141").Value = 17040
$ws.Range("M141").Value = -11860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -2340
$ws.Range("H122").Value = 2095.1177
$ws.Range("I122").Value = 1258.4286
$ws.Range("J122").Value = 5999.6665
$ws.Range("K122").Value = 3775.2858
$ws.Range("L122").Value = 17998.9995
$ws.Range("M122").Value = -1325.2858
$ws.Range("N122").Value = -22898.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1427.2858
$ws.Range("I86").Value = 1427.2858
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1427.2858
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -304.2858000000001
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1427.2858
$ws.Range("I89").Value = 1427.2858
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 7136.429
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1520.429
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1229.2222
$ws.Range("I5").Value = 302.16666
$ws.Range("K5").Value = 302.16666
$ws.Range("M5").Value = -190.16666
$ws.Range("H31").Value = 41872.2
$ws.Range("I31").Value = 2500
$ws.Range("K31").Value = 2500
$ws.Range("M31").Value = -2205
$ws.Range("H32").Value = 7829.3335
$ws.Range("I32").Value = 8244
$ws.Range("J32").Value = 7000
$ws.Range("K32").Value = 8244
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = -7928
$ws.Range("N32").Value = -7632
$ws.Range("H34").Value = 41872.2
$ws.Range("I34").Value = 2500
$ws.Range("K34").Value = 2500
$ws.Range("M34").Value = -2298
$ws.Range("H86").Value = 11936.526
$ws.Range("I86").Value = 12223.235
$ws.Range("J86").Value = 9499.5
$ws.Range("K86").Value = 12223.235
$ws.Range("L86").Value = 9499.5
$ws.Range("M86").Value = -11100.235
$ws.Range("N86").Value = -11745.5
$ws.Range("H89").Value = 11936.526
$ws.Range("I89").Value = 12223.235
$ws.Range("J89").Value = 9499.5
$ws.Range("K89").Value = 61116.175
$ws.Range("L89").Value = 47497.5
$ws.Range("M89").Value = -55500.175
$ws.Range("N89").Value = -58729.5
$ws.Range("H99").Value = 3417.5454
$ws.Range("I99").Value = 3259.3
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 3259.3
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -1761.3
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 3417.5454
$ws.Range("I126").Value = 3259.3
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9777.900000000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -7307.900000000001
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 4867.722
$ws.Range("I7").Value = 103.166664
$ws.Range("K7").Value = 309.499992
$ws.Range("M7").Value = -197.499992
$ws.Range("H9").Value = 180469.4
$ws.Range("I9").Value = 200147.5
$ws.Range("J9").Value = 167350.67
$ws.Range("K9").Value = 600442.5
$ws.Range("L9").Value = 502052.01
$ws.Range("M9").Value = -600218.5
$ws.Range("N9").Value = -502500.01
$ws.Range("H68").Value = 2127.4707
$ws.Range("J68").Value = 1402.4615
$ws.Range("L68").Value = 4207.3845
$ws.Range("N68").Value = -5829.3845
$ws.Range("H71").Value = 2127.4707
$ws.Range("J71").Value = 1402.4615
$ws.Range("L71").Value = 12622.1535
$ws.Range("N71").Value = -20734.1535
$ws.Range("H86").Value = 515.2593000000001
$ws.Range("J86").Value = 415.45456
$ws.Range("L86").Value = 1246.36368
$ws.Range("N86").Value = -3618.36368
$ws.Range("H89").Value = 515.2593000000001
$ws.Range("J89").Value = 415.45456
$ws.Range("L89").Value = 3739.09104
$ws.Range("N89").Value = -15595.09104
$ws.Range("H122").Value = 12662963
$ws.Range("I122").Value = 31145718
$ws.Range("J122").Value = 2581461
$ws.Range("K122").Value = 280311462
$ws.Range("L122").Value = 23233149
$ws.Range("M122").Value = -280309012
$ws.Range("N122").Value = -23238049

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5219.4
$ws.Range("I126").Value = 5935.769
$ws.Range("J126").Value = 3889
$ws.Range("K126").Value = 17807.307
$ws.Range("L126").Value = 11667
$ws.Range("M126").Value = -15337.307
$ws.Range("N126").Value = -16607

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3573.5862
$ws.Range("J22").Value = 2092.4
$ws.Range("L22").Value = 2092.4
$ws.Range("N22").Value = -2682.4
$ws.Range("H27").Value = 3573.5862
$ws.Range("J27").Value = 2092.4
$ws.Range("L27").Value = 2092.4
$ws.Range("N27").Value = -2306.4
$ws.Range("H33").Value = 19248.75
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H136").Value = 2673913
$ws.Range("I136").Value = 52774.75
$ws.Range("J136").Value = 3984482.2
$ws.Range("K136").Value = 158324.25
$ws.Range("L136").Value = 11953446.6
$ws.Range("M136").Value = -155774.25
$ws.Range("N136").Value = -11958546.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 219184.5
$ws.Range("I9").Value = 219184.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 219184.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -219044.5
$ws.Range("N9").ClearContents()
$ws.Range("H10").Value = 18999.5
$ws.Range("I10").Value = 18999.5
$ws.Range("K10").Value = 18999.5
$ws.Range("M10").Value = -18830.5
$ws.Range("H17").Value = 20000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 20000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -20344
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H122").Value = 2699
$ws.Range("I122").Value = 2932
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8796
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6346
$ws.Range("N122").Value = -10900

Write-Output "Applied scheduled profit data refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
